$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 53.463696812512026
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 46.05204738706685
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 60.76705279190513
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").Value = 46.609654277955656
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").Value = 39.785591828762811
$ws.Range("P9").Copy()
$ws.Range("Q9").PasteSpecial(-4122)
$ws.Range("Q9").Value = 53.693996785869842
$ws.Range("P10").Copy()
$ws.Range("Q10").PasteSpecial(-4122)
$ws.Range("Q10").Value = 49.132459991853935
$ws.Range("P11").Copy()
$ws.Range("Q11").PasteSpecial(-4122)
$ws.Range("Q11").Value = 42.132308166831223
$ws.Range("P12").Copy()
$ws.Range("Q12").PasteSpecial(-4122)
$ws.Range("Q12").Value = 56.225753650646354
$ws.Range("P13").Copy()
$ws.Range("Q13").PasteSpecial(-4122)
$ws.Range("Q13").Value = 28.457427087863305
$ws.Range("P14").Copy()
$ws.Range("Q14").PasteSpecial(-4122)
$ws.Range("Q14").Value = 20.524708126577082
$ws.Range("P15").Copy()
$ws.Range("Q15").PasteSpecial(-4122)
$ws.Range("Q15").Value = 36.325895173845353
$ws.Range("P16").Copy()
$ws.Range("Q16").PasteSpecial(-4122)
$ws.Range("Q16").Value = 37.816151622141014
$ws.Range("P17").Copy()
$ws.Range("Q17").PasteSpecial(-4122)
$ws.Range("Q17").Value = 29.032279844170926
$ws.Range("P18").Copy()
$ws.Range("Q18").PasteSpecial(-4122)
$ws.Range("Q18").Value = 46.928626462141906
$ws.Range("P19").Copy()
$ws.Range("Q19").PasteSpecial(-4122)
$ws.Range("Q19").Value = 51.38232216208695
$ws.Range("P20").Copy()
$ws.Range("Q20").PasteSpecial(-4122)
$ws.Range("Q20").Value = 45.862881450184311
$ws.Range("P21").Copy()
$ws.Range("Q21").PasteSpecial(-4122)
$ws.Range("Q21").Value = 57.0280888993139
$ws.Range("P22").Copy()
$ws.Range("Q22").PasteSpecial(-4122)
$ws.Range("Q22").Value = 44.951834666409091
$ws.Range("P23").Copy()
$ws.Range("Q23").PasteSpecial(-4122)
$ws.Range("Q23").Value = 38.216466887636237
$ws.Range("P24").Copy()
$ws.Range("Q24").PasteSpecial(-4122)
$ws.Range("Q24").Value = 51.83682668469686
$ws.Range("P25").Copy()
$ws.Range("Q25").PasteSpecial(-4122)
$ws.Range("Q25").Value = 82.176148450436926
$ws.Range("P26").Copy()
$ws.Range("Q26").PasteSpecial(-4122)
$ws.Range("Q26").Value = 66.965035434789911
$ws.Range("P27").Copy()
$ws.Range("Q27").PasteSpecial(-4122)
$ws.Range("Q27").Value = 96.931980629894966
$ws.Range("P28").Copy()
$ws.Range("Q28").PasteSpecial(-4122)
$ws.Range("Q28").Value = 56.391242440049062
$ws.Range("P29").Copy()
$ws.Range("Q29").PasteSpecial(-4122)
$ws.Range("Q29").Value = 50.844030930786069
$ws.Range("P30").Copy()
$ws.Range("Q30").PasteSpecial(-4122)
$ws.Range("Q30").Value = 61.300998533028128
$ws.Range("P31").Copy()
$ws.Range("Q31").PasteSpecial(-4122)
$ws.Range("Q31").Value = 54.829571415516767
$ws.Range("P32").Copy()
$ws.Range("Q32").PasteSpecial(-4122)
$ws.Range("Q32").Value = 58.407045187583961
$ws.Range("P33").Copy()
$ws.Range("Q33").PasteSpecial(-4122)
$ws.Range("Q33").Value = 51.452932817170577

$excel.CutCopyMode = $false
[void]$ws.Range("S34").Select()
